# Adiciona Caixa 42 com roupas íntimas (8 itens), prioridade Alta, Suíte.
# Insere as novas linhas em 342-349 na planilha "Catálogo de Mudança",
# empurrando as linhas existentes (antiga "Caixa 74" em diante) para baixo,
# e atualiza os totais na planilha "Resumo".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insere 8 novas linhas em branco a partir da linha 342
$ws.Range("A342:F349").EntireRow.Insert()

# Copia a formatação de uma linha existente de prioridade "Alta" (linha 2)
# para as novas linhas, garantindo os mesmos estilos (inclusive bordas).
$ws.Range("A2:F2").Copy()
$ws.Range("A342:F349").PasteSpecial(-4122)  # xlPasteFormats

$caixa = "Caixa 42"
$categoria = "Roupa íntima"
$comodo = "Suíte"
$prioridade = "Alta"

$descricoes = @(
    "Lingeries elaboradas",
    "Lingeries simples",
    "Calcinhas para período menstrual",
    "Cinta-liga",
    "Meia-calça",
    "Biquínis",
    "Maiôs",
    "Meias"
)

for ($i = 0; $i -lt $descricoes.Length; $i++) {
    $r = 342 + $i
    $ws.Cells.Item($r, 1).Value = $caixa
    $ws.Cells.Item($r, 2).Value = $categoria
    $ws.Cells.Item($r, 3).Value = $descricoes[$i]
    $ws.Cells.Item($r, 4).Value = $comodo
    $ws.Cells.Item($r, 5).Value = $prioridade
}

# Atualiza a planilha "Resumo": +8 itens no total e +8 itens de prioridade "Alta"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = 454
$ws2.Range("B6").Value = 104
